# Update CDA Logical model for ST.r2b
# Applies to the "Metadata" worksheet (Property/Value table):
#   - Version bumped 2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
#   - Date bumped to the new publication date/time
#   - A new "Jurisdiction" property row inserted right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update Version value (row 3) ---
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# --- Update Date value (row 8) ---
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" ---
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows (border/alignment/fill)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
